# Update the F-column "time_taken" timestamps on the "data" sheet.
# These represent a re-run of the panel query, each row's timestamp gets
# replaced with the value captured by the corresponding metadata refresh.
$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

$newTimestamps = @(
    "2021-10-05 14:33:58.879745",
    "2021-10-05 14:33:58.879755",
    "2021-10-05 14:33:58.879758",
    "2021-10-05 14:33:58.879761",
    "2021-10-05 14:33:58.879764",
    "2021-10-05 14:33:58.879766",
    "2021-10-05 14:33:58.879769",
    "2021-10-05 14:33:58.879772",
    "2021-10-05 14:33:58.879774",
    "2021-10-05 14:33:58.879777",
    "2021-10-05 14:33:58.879780",
    "2021-10-05 14:33:58.879782",
    "2021-10-05 14:33:58.879785",
    "2021-10-05 14:33:58.879787",
    "2021-10-05 14:33:58.879790",
    "2021-10-05 14:33:58.879793",
    "2021-10-05 14:33:58.879796",
    "2021-10-05 14:33:58.879798",
    "2021-10-05 14:33:58.879801",
    "2021-10-05 14:33:58.879803",
    "2021-10-05 14:33:58.879806",
    "2021-10-05 14:33:58.879809",
    "2021-10-05 14:33:58.879811",
    "2021-10-05 14:33:58.879814",
    "2021-10-05 14:33:58.879817",
    "2021-10-05 14:33:58.879820",
    "2021-10-05 14:33:58.879822",
    "2021-10-05 14:33:58.879825",
    "2021-10-05 14:33:58.879828",
    "2021-10-05 14:33:58.879830",
    "2021-10-05 14:33:58.879833",
    "2021-10-05 14:33:58.879835",
    "2021-10-05 14:33:58.879839",
    "2021-10-05 14:33:58.879841",
    "2021-10-05 14:33:58.879844",
    "2021-10-05 14:33:58.879846",
    "2021-10-05 14:33:58.879849",
    "2021-10-05 14:33:58.879852",
    "2021-10-05 14:33:58.879855",
    "2021-10-05 14:33:58.879857",
    "2021-10-05 14:33:58.879861",
    "2021-10-05 14:33:58.879863",
    "2021-10-05 14:33:58.879866",
    "2021-10-05 14:33:58.879869",
    "2021-10-05 14:33:58.879871",
    "2021-10-05 14:33:58.879874",
    "2021-10-05 14:33:58.879877",
    "2021-10-05 14:33:58.879879",
    "2021-10-05 14:33:58.879882",
    "2021-10-05 14:33:58.879884",
    "2021-10-05 14:33:58.879887",
    "2021-10-05 14:33:58.879889",
    "2021-10-05 14:33:58.879892",
    "2021-10-05 14:33:58.879895",
    "2021-10-05 14:33:58.879898",
    "2021-10-05 14:33:58.879901",
    "2021-10-05 14:33:58.879903",
    "2021-10-05 14:33:58.879906",
    "2021-10-05 14:33:58.879908",
    "2021-10-05 14:33:58.879911",
    "2021-10-05 14:33:58.879913",
    "2021-10-05 14:33:58.879916",
    "2021-10-05 14:33:58.879919",
    "2021-10-05 14:33:58.879921",
    "2021-10-05 14:33:58.879925",
    "2021-10-05 14:33:58.879929",
    "2021-10-05 14:33:58.879932",
    "2021-10-05 14:33:58.879934",
    "2021-10-05 14:33:58.879937",
    "2021-10-05 14:33:58.879940",
    "2021-10-05 14:33:58.879942",
    "2021-10-05 14:33:58.879945",
    "2021-10-05 14:33:58.879948",
    "2021-10-05 14:33:58.879950",
    "2021-10-05 14:33:58.879953",
    "2021-10-05 14:33:58.879955",
    "2021-10-05 14:33:58.879961",
    "2021-10-05 14:33:58.879964",
    "2021-10-05 14:33:58.879967",
    "2021-10-05 14:33:58.879970",
    "2021-10-05 14:33:58.879972",
    "2021-10-05 14:33:58.879975",
    "2021-10-05 14:33:58.879978",
    "2021-10-05 14:33:58.879981",
    "2021-10-05 14:33:58.879983",
    "2021-10-05 14:33:58.879986",
    "2021-10-05 14:33:58.879988",
    "2021-10-05 14:33:58.879991",
    "2021-10-05 14:33:58.879993",
    "2021-10-05 14:33:58.879996",
    "2021-10-05 14:33:58.879999",
    "2021-10-05 14:33:58.880001",
    "2021-10-05 14:33:58.880005",
    "2021-10-05 14:33:58.880008",
    "2021-10-05 14:33:58.880011",
    "2021-10-05 14:33:58.880014",
    "2021-10-05 14:33:58.880016",
    "2021-10-05 14:33:58.880019",
    "2021-10-05 14:33:58.880022",
    "2021-10-05 14:33:58.880024",
    "2021-10-05 14:33:58.880027",
    "2021-10-05 14:33:58.880030",
    "2021-10-05 14:33:58.880032",
    "2021-10-05 14:33:58.880035",
    "2021-10-05 14:33:58.880038",
    "2021-10-05 14:33:58.880040",
    "2021-10-05 14:33:58.880043",
    "2021-10-05 14:33:58.880045",
    "2021-10-05 14:33:58.880050",
    "2021-10-05 14:33:58.880054",
    "2021-10-05 14:33:58.880057",
    "2021-10-05 14:33:58.880059",
    "2021-10-05 14:33:58.880062",
    "2021-10-05 14:33:58.880065",
    "2021-10-05 14:33:58.880067",
    "2021-10-05 14:33:58.880070",
    "2021-10-05 14:33:58.880073",
    "2021-10-05 14:33:58.880075",
    "2021-10-05 14:33:58.880078",
    "2021-10-05 14:33:58.880081",
    "2021-10-05 14:33:58.880083",
    "2021-10-05 14:33:58.880086",
    "2021-10-05 14:33:58.880089",
    "2021-10-05 14:33:58.880091",
    "2021-10-05 14:33:58.880094",
    "2021-10-05 14:33:58.880097",
    "2021-10-05 14:33:58.880100",
    "2021-10-05 14:33:58.880102",
    "2021-10-05 14:33:58.880107",
    "2021-10-05 14:33:58.880110",
    "2021-10-05 14:33:58.880113",
    "2021-10-05 14:33:58.880116",
    "2021-10-05 14:33:58.880118"
)

for ($i = 0; $i -lt $newTimestamps.Length; $i++) {
    $row = $i + 2
    $dataSheet.Cells.Item($row, 6).Value = $newTimestamps[$i]
}

# Add a new "metadata" sheet right after "data" describing the panel query
# that produced this export.
$metaSheet = $wb.Worksheets.Add($null, $dataSheet)
$metaSheet.Name = "metadata"

$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

$metaSheet.Range("B1:G1").Font.Bold = $true
$metaSheet.Range("B1:G1").Borders.LineStyle = 1
$metaSheet.Range("B1:G1").HorizontalAlignment = -4108
$metaSheet.Range("B1:G1").VerticalAlignment = -4160

$metaSheet.Range("A2").Value = 0
$metaSheet.Range("A2").Font.Bold = $true
$metaSheet.Range("A2").Borders.LineStyle = 1
$metaSheet.Range("A2").HorizontalAlignment = -4108
$metaSheet.Range("A2").VerticalAlignment = -4160

$metaSheet.Range("B2").Value = "Hereditary Neuropathy - complex"
$metaSheet.Range("C2").Value = 3070
# Force "0.117" to be stored as text (matches the source export), then drop
# the now-superfluous explicit style so the cell stays unstyled like its
# siblings.
$metaSheet.Range("D2").NumberFormat = "@"
$metaSheet.Range("D2").Value = "0.117"
$metaSheet.Range("D2").Style = "Normal"
$metaSheet.Range("E2").Value = "2021-08-30T09:40:55.068399Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:33:58.875969"
$metaSheet.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/3070/?format=json"

# Keep "data" as the active sheet, matching the original workbook view.
$dataSheet.Activate()
